$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.702.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.06%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.635.72"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.75%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.05%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'217.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.55%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -0.76%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.01%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.73%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -0.74%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -0.37%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -0.20%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.865.22"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.67%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.644.36"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.60%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  -1.01%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.524"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.38%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'64.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.25%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'26.702.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.0₃0726"
$ws.Range("D18").Style = "Normal"
$ws.Range("E19").Value = "'  -0.01%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'210.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -3.25%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.57%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'2.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +2.94%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -1.49%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'9.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -2.72%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'145.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.26%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -0.30%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -2.01%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'7.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.72%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'15.55"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.09%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -2.15%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'1.19"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.68%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -0.31%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -1.10%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.276.61"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.15%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.52"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.36%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +0.08%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -2.00%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.531"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.52%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -1.79%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -0.12%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.801"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.63%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'2.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.38%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'1.775.84"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.58%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -3.64%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'60.73"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.55%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'91.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.54%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -1.78%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +0.64%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'7.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -2.44%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -0.87%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.406"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.44%  "
$ws.Range("E51").Style = "Normal"
